$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H8").Value = 175
$ws.Range("I8").Value = 175
$ws.Range("K8").Value = 525
$ws.Range("M8").Value = -386
$ws.Range("H40").Value = 2364.2856
$ws.Range("J40").Value = 2614.4285
$ws.Range("L40").Value = 2614.4285
$ws.Range("N40").Value = -2964.4285
$ws.Range("H106").Value = 4117945.5
$ws.Range("I106").Value = 4632318
$ws.Range("K106").Value = 4632318
$ws.Range("M106").Value = -4631687
$ws.Range("H118").Value = 294.33334
$ws.Range("I118").Value = 294.33334
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 883.0000200000001
$ws.Range("L118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = 773.9999799999999
$ws.Range("H137").Value = 32259532
$ws.Range("I137").Value = 55556708
$ws.Range("J137").Value = 1906
$ws.Range("K137").Value = 166670124
$ws.Range("L137").Value = 5718
$ws.Range("M137").Value = -166667574
$ws.Range("N137").Value = -10818

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20222.465
$ws.Range("I32").Value = 2306.6345
$ws.Range("J32").Value = 253128.25
$ws.Range("K32").Value = 2306.6345
$ws.Range("L32").Value = 253128.25
$ws.Range("M32").Value = -2019.6345
$ws.Range("N32").Value = -253702.25
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("H61").Value = 2419.2856
$ws.Range("I61").Value = 1484.08
$ws.Range("J61").Value = 4757.3
$ws.Range("K61").Value = 1484.08
$ws.Range("L61").Value = 4757.3
$ws.Range("M61").Value = -1272.08
$ws.Range("N61").Value = -5181.3
$ws.Range("H74").Value = 4398.35
$ws.Range("I74").Value = 1240.8572
$ws.Range("J74").Value = 26500.8
$ws.Range("K74").Value = 1240.8572
$ws.Range("L74").Value = 26500.8
$ws.Range("M74").Value = -366.8571999999999
$ws.Range("N74").Value = -28248.8
$ws.Range("H77").Value = 4398.35
$ws.Range("I77").Value = 1240.8572
$ws.Range("J77").Value = 26500.8
$ws.Range("K77").Value = 6204.286
$ws.Range("L77").Value = 132504
$ws.Range("M77").Value = -1836.286
$ws.Range("N77").Value = -141240
$ws.Range("H132").Value = 1798.3529
$ws.Range("I132").Value = 1535.5
$ws.Range("J132").Value = 3025
$ws.Range("K132").Value = 4606.5
$ws.Range("L132").Value = 9075
$ws.Range("M132").Value = -2076.5
$ws.Range("N132").Value = -14135
$ws.Range("H136").Value = 2419.2856
$ws.Range("I136").Value = 1484.08
$ws.Range("J136").Value = 4757.3
$ws.Range("K136").Value = 4452.24
$ws.Range("L136").Value = 14271.9
$ws.Range("M136").Value = -1902.24
$ws.Range("N136").Value = -19371.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 337.5
$ws.Range("I22").Value = 337.5
$ws.Range("K22").Value = 337.5
$ws.Range("M22").Value = -164.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 106.833336
$ws.Range("I7").Value = 110
$ws.Range("J7").Value = 105.25
$ws.Range("K7").Value = 110
$ws.Range("L7").Value = 105.25
$ws.Range("M7").Value = 3
$ws.Range("N7").Value = -331.25
$ws.Range("H22").Value = 596.875
$ws.Range("I22").Value = 225.66667
$ws.Range("J22").Value = 819.6
$ws.Range("K22").Value = 225.66667
$ws.Range("L22").Value = 819.6
$ws.Range("M22").Value = 124.33333
$ws.Range("N22").Value = -1519.6
$ws.Range("H31").Value = 1681.4706
$ws.Range("I31").Value = 873
$ws.Range("J31").Value = 3621.8
$ws.Range("K31").Value = 873
$ws.Range("L31").Value = 3621.8
$ws.Range("M31").Value = -578
$ws.Range("N31").Value = -4211.8
$ws.Range("H34").Value = 1681.4706
$ws.Range("I34").Value = 873
$ws.Range("J34").Value = 3621.8
$ws.Range("K34").Value = 873
$ws.Range("L34").Value = 3621.8
$ws.Range("M34").Value = -671
$ws.Range("N34").Value = -4025.8
$ws.Range("H39").Value = 60
$ws.Range("I39").Value = 60
$ws.Range("K39").Value = 60
$ws.Range("M39").Value = 331
$ws.Range("H49").Value = 60
$ws.Range("I49").Value = 60
$ws.Range("K49").Value = 60
$ws.Range("M49").Value = 122
$ws.Range("H58").Value = 1954.8478
$ws.Range("I58").Value = 744.4074000000001
$ws.Range("J58").Value = 3674.9473
$ws.Range("K58").Value = 744.4074000000001
$ws.Range("L58").Value = 3674.9473
$ws.Range("M58").Value = -541.4074000000001
$ws.Range("N58").Value = -4080.9473
$ws.Range("H132").Value = 2166.3274
$ws.Range("I132").Value = 1592
$ws.Range("J132").Value = 4463.636
$ws.Range("K132").Value = 4776
$ws.Range("L132").Value = 13390.908
$ws.Range("M132").Value = -2246
$ws.Range("N132").Value = -18450.908
$ws.Range("H136").Value = 1954.8478
$ws.Range("I136").Value = 744.4074000000001
$ws.Range("J136").Value = 3674.9473
$ws.Range("K136").Value = 2233.2222
$ws.Range("L136").Value = 11024.8419
$ws.Range("M136").Value = 316.7777999999998
$ws.Range("N136").Value = -16124.8419
$ws.Range("H138").Value = 52500
$ws.Range("J138").Value = 52500
$ws.Range("L138").Value = 52500
$ws.Range("N138").Value = -62780

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 4700
$ws.Range("I19").Value = 4700
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 14100
$ws.Range("L19").ClearContents()
$ws.Range("N19").Value = 0
$ws.Range("M19").Value = -13926
$ws.Range("H131").Value = 1800.9062
$ws.Range("I131").Value = 353.625
$ws.Range("J131").Value = 2283.3333
$ws.Range("K131").Value = 1060.875
$ws.Range("L131").Value = 6849.999899999999
$ws.Range("M131").Value = 3979.125
$ws.Range("N131").Value = -16929.9999
$ws.Range("H132").Value = 23810540
$ws.Range("I132").Value = 766.6667
$ws.Range("J132").Value = 41667868
$ws.Range("K132").Value = 6900.0003
$ws.Range("L132").Value = 375010812
$ws.Range("M132").Value = -4370.0003
$ws.Range("N132").Value = -375015872
$ws.Range("H133").Value = 5080
$ws.Range("I133").Value = 2114.2856
$ws.Range("K133").Value = 6342.8568
$ws.Range("M133").Value = -1282.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7416.1113
$ws.Range("I70").Value = 8349
$ws.Range("J70").Value = 6250
$ws.Range("K70").Value = 8349
$ws.Range("L70").Value = 6250
$ws.Range("M70").Value = -8079
$ws.Range("N70").Value = -6790
$ws.Range("H73").Value = 7416.1113
$ws.Range("I73").Value = 8349
$ws.Range("J73").Value = 6250
$ws.Range("K73").Value = 8349
$ws.Range("L73").Value = 6250
$ws.Range("M73").Value = -7413
$ws.Range("N73").Value = -8122
$ws.Range("H132").Value = 3113.2563
$ws.Range("I132").Value = 2876.6177
$ws.Range("J132").Value = 4722.4
$ws.Range("K132").Value = 8629.8531
$ws.Range("L132").Value = 14167.2
$ws.Range("M132").Value = -6099.8531
$ws.Range("N132").Value = -19227.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16274

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 85114.836
$ws.Range("J122").Value = 2100
$ws.Range("L122").Value = 6300
$ws.Range("N122").Value = -11200

